$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HomeResource")

# Remove outdated/untranslated "detail" (column F) entries that no longer apply
$ws.Range("F3").ClearContents()
$ws.Range("F4").ClearContents()
$ws.Range("F7").ClearContents()

# Update the active selection to match the latest edit position
$ws.Range("F4").Select()
